# Remove the "20項 (20 things about Japan)" post row (row 227) from the
# posts sheet. All rows below shift up by one, which Excel's row-delete
# handles natively (matches the source diff: old row 227 is gone, old rows
# 228-251 become new rows 227-250, dimension shrinks from C251 to C250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(227).Delete()
